$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.5926966292134831
$ws1.Range("C2").Value = 0.5521601685985248
$ws1.Range("D2").Value = 0.9812734082397003
$ws1.Range("E2").Value = 0.7066756574511126
$ws1.Range("F2").Value = 0.8492706645056726
$ws1.Range("G2").Value = 0.9527939016714455
$ws1.Range("H2").Value = 0.7908969125671562
$ws1.Range("I2").Value = 524
$ws1.Range("J2").Value = 425
$ws1.Range("K2").Value = 109
$ws1.Range("L2").Value = 10

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

$ws2.Range("B2").Value = 0.9159663865546218
$ws2.Range("C2").Value = 0.2041198501872659
$ws2.Range("D2").Value = 0.333843797856049

$ws2.Range("B3").Value = 0.5521601685985248
$ws2.Range("C3").Value = 0.9812734082397003
$ws2.Range("D3").Value = 0.7066756574511126

$ws2.Range("B4").Value = 0.5926966292134831
$ws2.Range("C4").Value = 0.5926966292134831
$ws2.Range("D4").Value = 0.5926966292134831
$ws2.Range("E4").Value = 0.5926966292134831

$ws2.Range("B5").Value = 0.7340632775765733
$ws2.Range("C5").Value = 0.5926966292134831
$ws2.Range("D5").Value = 0.5202597276535808

$ws2.Range("B6").Value = 0.7340632775765733
$ws2.Range("C6").Value = 0.5926966292134831
$ws2.Range("D6").Value = 0.5202597276535808

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

$ws3.Range("B2").Value = 109
$ws3.Range("C2").Value = 425

$ws3.Range("B3").Value = 10
$ws3.Range("C3").Value = 524
